# Add team record (Wins/Losses/Ties) columns to the data sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 43

# Header row (styled like the other headers in row 1).
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

# Data rows: every row gets the same team record.
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 30).Value = 73
    $ws.Cells.Item($r, 31).Value = 89
    $ws.Cells.Item($r, 32).Value = 0
}
